$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, shifting existing rows 12:85 down to 13:86.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new data record.
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(12, 3).Value = 'La Araucanía'
$ws.Cells.Item(12, 4).Value = 44901
$ws.Cells.Item(12, 5).Value = 9
$ws.Cells.Item(12, 6).Value = 300000000
$ws.Cells.Item(12, 7).Value = 'Espárragos'
$ws.Cells.Item(12, 8).Value = 'Sin especificar'
$ws.Cells.Item(12, 9).Value = 'Primera'
$ws.Cells.Item(12, 10).Value = 250
$ws.Cells.Item(12, 11).Value = 1500
$ws.Cells.Item(12, 12).Value = 1500
$ws.Cells.Item(12, 13).Value = 1500
$ws.Cells.Item(12, 14).Value = '$/kilo'
$ws.Cells.Item(12, 15).Value = 'Provincia de Linares'
$ws.Cells.Item(12, 16).Value = 1500
$ws.Cells.Item(12, 17).Value = 1
$ws.Cells.Item(12, 18).Value = 'Hortaliza'
